$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteAll = -4104
$xlPasteFormats = -4122

# --- Reshuffle existing cell contents (value+format) into their new homes. ---
# Order matters: every source cell is read before anything overwrites it.
$ws.Range("B24").Copy() | Out-Null
$ws.Range("B23").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B24").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C24").Copy() | Out-Null
$ws.Range("C23").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C25").Copy() | Out-Null
$ws.Range("C24").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A17").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A18").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A19").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A20").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A21").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A23").Copy() | Out-Null
$ws.Range("A22").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B18").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C18").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B14").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C14").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C15").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A13").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A14").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A15").PasteSpecial($xlPasteAll) | Out-Null

# --- New text that did not exist anywhere before ("Semestral") ---
$ws.Range("B13").Value = "Semestral"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C13").Value = "Semestral"
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C13").PasteSpecial($xlPasteFormats) | Out-Null

# --- Clear cells that must become empty ---
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# --- Row heights to match the new layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30

# Rows 17 and 22 revert to the default (no explicit custom height).
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(22).AutoFit()

# --- Remove the now-obsolete last row (25) ---
$ws.Rows.Item(25).Delete()
